# Update data refresh timestamp on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 12:26 PM"

# Update the "1 Year" column (F) on the Industry Analysis sheet for rows 2-76.
$industry = $wb.Worksheets.Item("Industry Analysis")

$newValues = @(
    18.476,
    -7.7404,
    30.7972,
    -50.2266,
    61.9649,
    -9.1713,
    -3.556,
    38.3509,
    -6.2497,
    52.6723,
    -6.932,
    17.5662,
    -35.5106,
    0.6286,
    -3.1514,
    -20.6354,
    -0.0175,
    -26.9255,
    44.703,
    10.0506,
    84.6016,
    -54.4868,
    -12.8122,
    -9.182700000000001,
    5.9529,
    -33.2998,
    -20.4441,
    -17.1514,
    24.527,
    57.6193,
    -1.527,
    -5.2378,
    27.4054,
    6.7961,
    -5.6683,
    1.4178,
    -22.4272,
    12.3741,
    -5.138,
    -0.1825,
    23.2483,
    14.456,
    -11.1739,
    27.112,
    -5.6252,
    -36.5148,
    -27.8397,
    -25.4424,
    -49.1173,
    -51.065,
    -35.4517,
    -11.9879,
    -3.0992,
    -15.3441,
    -25.937,
    -29.1486,
    -6.4093,
    -23.3046,
    -11.2657,
    -9.777699999999999,
    -16.0561,
    -9.932499999999999,
    51.8767,
    -43.5191,
    13.7315,
    12.6111,
    31.7532,
    -19.9577,
    -12.9642,
    13.2432,
    2.8232,
    -9.179,
    -14.2931,
    28.3699,
    45.5868
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $industry.Cells.Item($row, 6).Value = $newValues[$i]
}
